# Add data for 2022-05-21: refresh the "through May 12" rolling window to
# "through May 13" and bump/insert the corresponding carjacking counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the header label for the rolling-window column.
$ws.Name = "Through 2022-05-13"
$ws.Range("B1").Value = "May 2022 (through May 13)"

# Updated / new counts in the "May 2022 (through May 13)" data block.
$ws.Range("G3").Value = 6      # Austin
$ws.Range("B4").Value = 2      # Humboldt Park
$ws.Range("V4").Value = 2      # Humboldt Park
$ws.Range("G5").Value = 2      # Garfield Park
$ws.Range("L5").Value = 2      # Garfield Park
$ws.Range("AK5").Value = 1     # Garfield Park
$ws.Range("V6").Value = 1      # Chicago Lawn
$ws.Range("B7").Value = 1      # North Lawndale
$ws.Range("B8").Value = 3      # South Shore
$ws.Range("AA8").Value = 1     # South Shore
$ws.Range("V10").Value = 1     # Belmont Cragin
$ws.Range("AK14").Value = 1    # Lincoln Park
$ws.Range("G16").Value = 1     # Little Italy, UIC
$ws.Range("AA17").Value = 1    # South Chicago
$ws.Range("G21").Value = 1     # Chatham
$ws.Range("B23").Value = 3     # Grand Crossing
$ws.Range("L29").Value = 1     # West Pullman
$ws.Range("AK52").Value = 1    # Beverly
$ws.Range("G85").Value = 1     # Rogers Park
$ws.Range("B91").Value = 2     # Washington Park
